# "tiny changes for default methods"
#
# Slide 2, shape "CustomShape 2" has three numbered-list paragraphs:
#   1) Create a new Interface with one abstract method and one default method
#   2) Create a new class implementing your interface
#   3) Call both your implemented method and the default method
#
# The diff tweaks the wording of paragraphs (1) and (3) while leaving each
# run's character formatting (font color / spell-check "err" flag) alone
# wherever the run itself is otherwise untouched. We locate the exact
# substrings with TextRange.Find(...) *before* making any edits (so every
# offset is resolved against the original text), then rewrite them via
# TextRange.Characters(start,length).Text, applying the edits from the
# end of the text frame towards the start so that earlier offsets are
# never invalidated by a preceding (in textual order) edit changing the
# overall text length.
#
# Note: TextRange.Find(what, after) returns the first match that starts
# strictly AFTER the given position, so chained Find() calls use
# (previous.Start + previous.Length - 1) as the next "after" argument.

$p  = $ppt.ActivePresentation
$s  = $p.Slides.Item(2)
$sh = $s.Shapes.Item(2)
$tr = $sh.TextFrame.TextRange

# --- locate everything up-front, against the untouched original text ---

# Paragraph 1
$fAbstract = $tr.Find(" abstract ")
$fMethod1  = $tr.Find("method", $fAbstract.Start + $fAbstract.Length - 1)
$fAnd      = $tr.Find(" and ")
$fOne2     = $tr.Find("one", $fAnd.Start + $fAnd.Length - 1)
$fDefault1 = $tr.Find(" default ")

# Paragraph 3
$fCall        = $tr.Find("Call ")
$fBoth        = $tr.Find("both", $fCall.Start + $fCall.Length - 1)
$fYour1       = $tr.Find("your", $fBoth.Start + $fBoth.Length - 1)
$fImplemented = $tr.Find("implemented", $fYour1.Start + $fYour1.Length - 1)
$fMethod2     = $tr.Find("method", $fImplemented.Start + $fImplemented.Length - 1)
$fAndDefault  = $tr.Find(" and the default ", $fMethod2.Start + $fMethod2.Length - 1)

# --- apply the edits, highest Start first ---

# Paragraph 3, back to front
# " and the default " -> "default "
$tr.Characters($fAndDefault.Start, $fAndDefault.Length).Text = "default "

# "implemented" + " " + "method" -> removed entirely
$removeLen = ($fMethod2.Start + $fMethod2.Length) - $fImplemented.Start
$tr.Characters($fImplemented.Start, $removeLen).Text = ""

# second "your" -> removed entirely
$tr.Characters($fYour1.Start, $fYour1.Length).Text = ""

# "both" -> "your"
$tr.Characters($fBoth.Start, $fBoth.Length).Text = "your"

# Paragraph 1, back to front
# " default " -> "default "
$tr.Characters($fDefault1.Start, $fDefault1.Length).Text = "default "

# " and " + "one" -> " " (the two runs collapse into one, keeping the
# formatting of the " and " run)
$mergeLen = ($fOne2.Start + $fOne2.Length) - $fAnd.Start
$tr.Characters($fAnd.Start, $mergeLen).Text = " "

# "method" -> "one"
$tr.Characters($fMethod1.Start, $fMethod1.Length).Text = "one"

# " abstract " -> " "
$tr.Characters($fAbstract.Start, $fAbstract.Length).Text = " "

# --- finally, insert the new run (a single space) right after "Call " ---
# (done last, after every other rewrite above; earlier edits to paragraph 1
# shift the absolute character offsets of paragraph 3, so re-find "Call "
# fresh against the current text instead of reusing the stale $fCall.Start)
$trFinal = $sh.TextFrame.TextRange
$fCallFinal = $trFinal.Find("Call ")
$lastCharOfCall = $trFinal.Characters($fCallFinal.Start + $fCallFinal.Length - 1, 1)
$lastCharOfCall.InsertAfter(" ")
